# Updated code on timestamp: 01-11-2021 - 11:09:29.15
#
# This script reproduces the authoring change: a batch "clear contents" pass
# over a handful of scratch/helper columns (F/G) on the Stack, Queue,
# Linked List and Generic Tree sheets, plus the corresponding change of
# active sheet/selection that Excel records when a user does this kind of
# cleanup (ending up with the "Stack" sheet active/selected).

$wb = $excel.ActiveWorkbook

# ---- Queue sheet: clear helper column F (rows 4-10) ----
$wsQueue = $wb.Worksheets.Item("Queue")
$wsQueue.Range("F4:F10").ClearContents()
[void]$wsQueue.Range("F4:F10").Select()

# ---- Linked List sheet: clear helper column G (rows 7-49) ----
$wsLinkedList = $wb.Worksheets.Item("Linked List")
$wsLinkedList.Range("G7:G49").ClearContents()
[void]$wsLinkedList.Range("G1:G1048576").Select()

# ---- Generic Tree sheet: clear helper column F (rows 3-32) ----
$wsGenericTree = $wb.Worksheets.Item("Generic Tree")
$wsGenericTree.Range("F3:F32").ClearContents()
[void]$wsGenericTree.Range("F1:F1048576").Select()

# ---- Stack sheet: clear helper column F (rows 5-21), then make it the
# active sheet with F5:F21 selected (this is the last sheet touched, so it
# becomes the workbook's active tab / tabSelected sheet) ----
$wsStack = $wb.Worksheets.Item("Stack")
$wsStack.Range("F5:F21").ClearContents()
[void]$wsStack.Activate()
[void]$wsStack.Range("F5:F21").Select()
